$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
  2  = -2
  3  = -1
  4  = 6
  5  = -3
  6  = 3
  7  = -2
  10 = -5
  11 = 8
  12 = 2
  13 = -1
  14 = -3
  15 = -2
  16 = 0
  17 = 1
  19 = 0
  20 = -3
  21 = -4
  22 = -3
  23 = -7
  24 = 4
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
